# Apply the "Updated cryptos list" refresh: new prices / hourly volume
# deltas for every coin row, plus three coins (Aptos, InternetComputer,
# Hedera, FraxShare) shifting down one rank and Quant being replaced by
# WEMIXTOKEN in the last row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.798.82'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '2.116.18'
$ws.Range("E3").Value = '  +6.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.26'
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5333'
$ws.Range("E7").Value = '  +4.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4406'
$ws.Range("E8").Value = '  +6.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09013'
$ws.Range("E9").Value = '  +3.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.78'
$ws.Range("E10").Value = '  +12.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.179'
$ws.Range("E11").Value = '  +4.39%  '
$ws.Range("E12").Value = '  +3.24%  '
$ws.Range("D13").Value = '2.114.29'
$ws.Range("E13").Value = '  +6.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.770'
$ws.Range("E14").Value = '  +4.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.822'
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.68'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06676'
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.340'
$ws.Range("E22").Value = '  +4.48%  '
$ws.Range("D23").Value = '30.868.33'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.33'
$ws.Range("E24").Value = '  +6.50%  '
$ws.Range("D25").Value = '2.360.36'
$ws.Range("E25").Value = '  +6.10%  '
$ws.Range("E26").Value = '  +3.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.74'
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.589'
$ws.Range("E28").Value = '  +9.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.37'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.39'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.189'
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.228'
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.010'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.562'
$ws.Range("E35").Value = '  +18.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02602'
$ws.Range("E36").Value = '  +4.95%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.544'
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06772'
$ws.Range("E38").Value = '  +3.90%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.571'
$ws.Range("E39").Value = '  +7.60%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.86'
$ws.Range("E40").Value = '  +9.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2320'
$ws.Range("E41").Value = '  +5.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6851'
$ws.Range("E42").Value = '  +4.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.248'
$ws.Range("E43").Value = '  +1.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6467'
$ws.Range("E44").Value = '  +5.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.16'
$ws.Range("E45").Value = '  +3.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("E47").Value = '  +1.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.670'
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.266'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.03'
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.198'
$ws.Range("E51").Value = '  +8.98%  '
